# Account balance, column additional effort removed
#
# The "Additional Effort [h]" column (C) is folded into the "Effort [h]"
# column (B) - any hours recorded as additional effort are added on top of
# the regular effort for that day - and the now-empty/obsolete column C is
# removed, shifting the former "comment" column D left into C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fold the additional-effort hours (column C) into the effort hours
# (column B) for every row that had a value in column C.
$ws.Range("B7").Value  = 4      # 3    + 1
$ws.Range("B8").Value  = 6      # 5    + 1
$ws.Range("B10").Value = 2.5    # 1.5  + 1
$ws.Range("B11").Value = 1      # B11 was empty, C11 held 1
$ws.Range("B18").Value = 4      # 2    + 2
$ws.Range("B31").Value = 3.75   # 2.75 + 1

# Now that its values have been merged into column B, delete the
# "Additional Effort [h]" column outright; this shifts the comment
# column (D) left into column C and drops the orphaned shared string.
$ws.Columns.Item(3).Delete()
